# Append the new resale-number row (row 44) to the CityResaleNum sheet,
# matching the 2025-01-26 23:14:55 snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 44

# Text columns (A-D): force Text format while assigning so values like the
# date string and the leading-zero week "04" aren't auto-converted into a
# date serial / number, then clear the (now unneeded) formatting so the
# cells end up styled the same as the rest of the data (no explicit style).
$textRange = $ws.Range("A" + $row + ":D" + $row)
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-01-26"
$ws.Cells.Item($row, 2).Value = "23:14:55"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 4).Value = "04"

$textRange.ClearFormats()

# Numeric columns (E-T)
$ws.Cells.Item($row, 5).Value = 126099
$ws.Cells.Item($row, 6).Value = 142011
$ws.Cells.Item($row, 7).Value = 168022
$ws.Cells.Item($row, 8).Value = 158479
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142567
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191523
$ws.Cells.Item($row, 14).Value = 115616
$ws.Cells.Item($row, 15).Value = 45494
$ws.Cells.Item($row, 16).Value = 28368
$ws.Cells.Item($row, 17).Value = 64913
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 45915
$ws.Cells.Item($row, 20).Value = -1
